$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F ("想去人数" / want-to-go count) updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 683
$wsExhibit.Range("F4").Value = 1992
$wsExhibit.Range("F5").Value = 5821
$wsExhibit.Range("F8").Value = 3289
$wsExhibit.Range("F12").Value = 4581
$wsExhibit.Range("F16").Value = 6
$wsExhibit.Range("F17").Value = 52
$wsExhibit.Range("F22").Value = 310
$wsExhibit.Range("F32").Value = 209
$wsExhibit.Range("F33").Value = 412
$wsExhibit.Range("F37").Value = 2263
$wsExhibit.Range("F38").Value = 1057
$wsExhibit.Range("F42").Value = 387

# Sheet "演出" (Performance) - column F update
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F11").Value = 159

# Sheet "全部类型" (All types) - column F updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 683
$wsAll.Range("F5").Value = 1992
$wsAll.Range("F6").Value = 5821
$wsAll.Range("F9").Value = 3289
$wsAll.Range("F12").Value = 4581
$wsAll.Range("F14").Value = 6
$wsAll.Range("F16").Value = 52
$wsAll.Range("F23").Value = 310
$wsAll.Range("F31").Value = 209
$wsAll.Range("F34").Value = 2263
$wsAll.Range("F35").Value = 1057
$wsAll.Range("F41").Value = 387
